$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Paragraph A (ilvl=1): intro sentence about running python without
# installing it, with a hyperlink to https://replit.com
# ---------------------------------------------------------------------
$last = $d.Paragraphs.Last
$last.Range.InsertParagraphAfter()
$paraA = $d.Paragraphs.Item($d.Paragraphs.Count)
$paraA.Range.ListFormat.ListLevelNumber = 2

$rng = $paraA.Range
$rng.Collapse(0)
$fullText = " Caso não possa instalar o python na sua máquina por qualquer motivo que seja, você pode utilizá-lo no navegador através do site: https://replit.com. Lá você também poderá escrever linhas de código."
$rng.InsertAfter($fullText)

$linkStart = $paraA.Range.Start + $fullText.IndexOf("https://replit.com")
$linkEnd = $linkStart + [string]"https://replit.com".Length
$linkRange = $d.Range($linkStart, $linkEnd)
$d.Hyperlinks.Add($linkRange, "https://replit.com") | Out-Null
$newLink = $d.Hyperlinks.Item($d.Hyperlinks.Count)
$newLink.Range.Font.Name = "Times New Roman"
$newLink.Range.Font.Size = 12

# ---------------------------------------------------------------------
# Paragraph B (ilvl=2): the site serves several languages
# ---------------------------------------------------------------------
$paraA.Range.InsertParagraphAfter()
$paraB = $d.Paragraphs.Item($d.Paragraphs.Count)
$paraB.Range.ListFormat.ListLevelNumber = 3
$rng = $paraB.Range
$rng.Collapse(0)
$rng.InsertAfter("Esse site serve para diversas linguagens, não somente pytho.")

# ---------------------------------------------------------------------
# Paragraph C (ilvl=2): Python3 is the first option after "Start coding"
# ---------------------------------------------------------------------
$paraB.Range.InsertParagraphAfter()
$paraC = $d.Paragraphs.Item($d.Paragraphs.Count)
$paraC.Range.ListFormat.ListLevelNumber = 3
$rng = $paraC.Range
$rng.Collapse(0)
$rng.InsertAfter("Python3 vem como primeira opção após clicar em " + [char]0x201C + "Start coding" + [char]0x201D + ".")

# ---------------------------------------------------------------------
# Paragraph D (ilvl=2): one of the best sites to code
# ---------------------------------------------------------------------
$paraC.Range.InsertParagraphAfter()
$paraD = $d.Paragraphs.Item($d.Paragraphs.Count)
$paraD.Range.ListFormat.ListLevelNumber = 3
$rng = $paraD.Range
$rng.Collapse(0)
$rng.InsertAfter("Este é um dos melhores sites para codar.")

# ---------------------------------------------------------------------
# Paragraph E (ilvl=1): trailing blank paragraph
# ---------------------------------------------------------------------
$paraD.Range.InsertParagraphAfter()
$paraE = $d.Paragraphs.Item($d.Paragraphs.Count)
$paraE.Range.ListFormat.ListLevelNumber = 2
$rng = $paraE.Range
$rng.Collapse(0)
$rng.InsertAfter(" ")

Write-Output "ok"
